# Generate Report for Handoff
#
# Status text moves from "Handed back: in sync with en-US" to
# "Ready for handoff", and the associated handoff timestamps are bumped
# to reflect the freshly generated report. The "date-ish" status/date
# columns are also narrowed now that the new status text is shorter.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status + generate date ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-17 16:58:18"

# --- zh-cn detail sheet: status + latest handoff datetime ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-17 16:58:13"

# --- de-de detail sheet: status + latest handoff datetime ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-17 16:58:18"

# --- Narrow the status/date columns to fit the new, shorter content ---
$overview.Range("E1").ColumnWidth = 16.3
$overview.Range("F1").ColumnWidth = 16.3

$zhcn.Range("C1").ColumnWidth = 16.3

$dede.Range("C1").ColumnWidth = 16.3
